$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cells in column D may contain plain decimal-looking text (e.g. "272.98")
# which Excel would otherwise auto-convert to a number on assignment.
# Force text format, assign the literal string, then restore the default
# "Normal" style so the cell style index is unaffected.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '43.708.82'
$ws.Range('E2').Value = '  -0.83%  '
Set-TextValue 'D3' '2.229.36'
$ws.Range('E3').Value = '  -1.09%  '
$ws.Range('E4').Value = '  -0.06%  '
Set-TextValue 'D5' '272.98'
$ws.Range('E5').Value = '  +5.63%  '
Set-TextValue 'D6' '87.08'
$ws.Range('E6').Value = '  +9.77%  '
$ws.Range('E7').Value = '  -0.93%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -0.75%  '
Set-TextValue 'D10' '44.88'
$ws.Range('E10').Value = '  +3.69%  '
Set-TextValue 'D11' '0.0919'
$ws.Range('E11').Value = '  -1.17%  '
Set-TextValue 'D12' '7.69'
$ws.Range('E12').Value = '  +8.13%  '
$ws.Range('E13').Value = '  +1.12%  '
Set-TextValue 'D14' '2.565.41'
$ws.Range('E14').Value = '  -1.00%  '
Set-TextValue 'D15' '14.87'
$ws.Range('E15').Value = '  +0.57%  '
Set-TextValue 'D16' '2.216.70'
$ws.Range('E16').Value = '  -0.58%  '
Set-TextValue 'D17' '0.790'
$ws.Range('E17').Value = '  -0.83%  '
Set-TextValue 'D18' '43.652.46'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('E19').Value = '  -1.67%  '
$ws.Range('E20').Value = '  -1.99%  '
Set-TextValue 'D21' '5.95'
$ws.Range('E22').Value = '  -0.59%  '
Set-TextValue 'D23' '232.24'
$ws.Range('E23').Value = '  -1.23%  '
Set-TextValue 'D24' '8.67'
$ws.Range('E24').Value = '  -9.06%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('E26').Value = '  +13.20%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D27' '10.74'
$ws.Range('E27').Value = '  -0.65%  '
$ws.Range('B28').Value = 'WEMIXToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D28' '3.53'
$ws.Range('E28').Value = '  +4.45%  '
$ws.Range('E29').Value = '  +5.31%  '
Set-TextValue 'D30' '39.16'
$ws.Range('E30').Value = '  -7.57%  '
Set-TextValue 'D31' '173.21'
$ws.Range('E31').Value = '  -0.21%  '
Set-TextValue 'D32' '0.0901'
$ws.Range('E32').Value = '  +2.63%  '
Set-TextValue 'D33' '20.37'
$ws.Range('E33').Value = '  -1.64%  '
Set-TextValue 'D34' '5.34'
$ws.Range('E34').Value = '  +0.15%  '
$ws.Range('E35').Value = '  -0.23%  '
$ws.Range('E36').Value = '  -3.04%  '
Set-TextValue 'D37' '0.0351'
$ws.Range('E37').Value = '  -3.23%  '
Set-TextValue 'D38' '4.30'
$ws.Range('E38').Value = '  -4.89%  '
Set-TextValue 'D39' '3.35'
$ws.Range('E39').Value = '  +16.23%  '
Set-TextValue 'D40' '2.20'
$ws.Range('E40').Value = '  +2.18%  '
Set-TextValue 'D41' '12.38'
$ws.Range('E41').Value = '  -8.17%  '
Set-TextValue 'D42' '63.52'
$ws.Range('E42').Value = '  +1.29%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('E44').Value = '  -0.85%  '
Set-TextValue 'D45' '8.47'
$ws.Range('E45').Value = '  -1.34%  '
Set-TextValue 'D46' '0.0987'
$ws.Range('E46').Value = '  -0.27%  '
Set-TextValue 'D47' '99.89'
$ws.Range('E47').Value = '  -5.08%  '
$ws.Range('E48').Value = '  +3.16%  '
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('B50').Value = 'WOONetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
Set-TextValue 'D50' '0.428'
$ws.Range('E50').Value = '  -9.53%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D51' '1.49'
$ws.Range('E51').Value = '  -1.87%  '
